$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 277.4
$ws.Range("C3").Value = 331.8
$ws.Range("B4").Value = 122.5
$ws.Range("C4").Value = 220
$ws.Range("C5").Value = 153.5
$ws.Range("C6").Value = -114.6
$ws.Range("C7").Value = -235.3
$ws.Range("C8").Value = -184
$ws.Range("C9").Value = -140
$ws.Range("C10").Value = -380.6
$ws.Range("C11").Value = -269.3
$ws.Range("C12").Value = -228
$ws.Range("C13").Value = -356.3
$ws.Range("C14").Value = 135.8
$ws.Range("C15").Value = 583.9
$ws.Range("C16").Value = 572.2
$ws.Range("C17").Value = 55.9
$ws.Range("C18").Value = -120.8
$ws.Range("C19").Value = 223.3
$ws.Range("C20").Value = -72.90000000000001
$ws.Range("C21").Value = 11.3
$ws.Range("C22").Value = 414.1
$ws.Range("C23").Value = 464.2
$ws.Range("C24").Value = 357.3
$ws.Range("C25").Value = 154.3
